$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Revise existing covid_deaths (column C) counts for a number of
#        previously-reported rows (upward revisions of historical data). ---
$cUpdates = @{
    992  = 35
    1044 = 30
    1058 = 22
    1077 = 2
    1084 = 23
    1085 = 48
    1095 = 13
    1096 = 38
    1099 = 10
    1104 = 22
    1105 = 34
    1107 = 9
    1109 = 34
    1113 = 30
    1115 = 13
    1116 = 22
    1117 = 31
}

foreach ($r in $cUpdates.Keys) {
    $ws.Cells.Item($r, 3).Value = $cUpdates[$r]
}

# --- 2) The last 3 rows (1118-1120), which used to be the final rows for
#        date 44187 (2020-12-22), are rewritten with corrected agegrp/
#        covid_deaths values. ---
$ws.Cells.Item(1118, 2).Value = "0-19"
$ws.Cells.Item(1118, 3).Value = 1

$ws.Cells.Item(1119, 2).Value = "40-49"
$ws.Cells.Item(1119, 3).Value = 1

$ws.Cells.Item(1120, 2).Value = "50-59"
$ws.Cells.Item(1120, 3).Value = 1

# --- 3) Append 6 brand-new rows (1121-1126): three more agegrp rows to
#        finish out date 44187, then three rows for the new date 44188
#        (2020-12-23). Column A reuses the same date number format (style)
#        as the other date cells in the sheet. ---
$dateNumberFormat = $ws.Cells.Item(1117, 1).NumberFormat()

$newRows = @(
    @{ Row = 1121; Date = 44187; AgeGrp = "60-69"; Deaths = 9 },
    @{ Row = 1122; Date = 44187; AgeGrp = "70-79"; Deaths = 16 },
    @{ Row = 1123; Date = 44187; AgeGrp = "80+";   Deaths = 18 },
    @{ Row = 1124; Date = 44188; AgeGrp = "60-69"; Deaths = 1 },
    @{ Row = 1125; Date = 44188; AgeGrp = "70-79"; Deaths = 12 },
    @{ Row = 1126; Date = 44188; AgeGrp = "80+";   Deaths = 9 }
)

foreach ($nr in $newRows) {
    $r = $nr.Row
    $ws.Cells.Item($r, 1).Value = $nr.Date
    $ws.Cells.Item($r, 1).NumberFormat = $dateNumberFormat
    $ws.Cells.Item($r, 2).Value = $nr.AgeGrp
    $ws.Cells.Item($r, 3).Value = $nr.Deaths
}
